$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws1.Range("B6").Value = "C:\Users\Marina Cernat\Documents\GitHub\rpa-testing\TestingTool_v3"
$ws1.Range("B8").Value = "C:\Users\Marina Cernat\Documents\GitHub\rpa-testing\TestingTool_v3\first.py"
$ws1.Range("B10").Value = "C:\Users\Marina Cernat\Documents\GitHub\rpa-testing\TestingTool_v3\Applications\C#Models\SimpleBankLoanCSharp"
$ws1.Range("B4").Value = "C:\Users\Marina Cernat\AppData\Local\Programs\Python\Python38"
$ws1.Range("B4").WrapText = $False

$ws1.Columns.Item(2).ColumnWidth = 110.3

$null = $ws1.Activate()
$null = $ws1.Range("A1").Select()
$null = $ws1.Range("B16").Select()
